# edit.ps1 - applies:
#   1) Date placeholder text "20-11-2023" -> "15/08/24" (slide master + every layout)
#   2) Slide background C00000 -> 821918 (all 10 slides)
#   3) The two side "Rectangle" accent shapes: bg2/lumMod25% -> solid B08A4B (all 10 slides)
#   4) Slide 6's permissions table: every cell fill C00000 -> 821918, and a tiny
#      table-frame height correction picked up by PowerPoint's relayout.

function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$ACCENT_BGR = ToRGB 0x82 0x19 0x18   # srgbClr 821918
$GOLD_BGR   = ToRGB 0xB0 0x8A 0x4B   # srgbClr B08A4B
$NEW_DATE   = "15/08/24"

$p = $ppt.ActivePresentation

# --- 1) Refresh the cached "Date Placeholder" text everywhere it appears ---
$m = $p.SlideMaster

for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $NEW_DATE
    }
}

for ($i = 1; $i -le $m.CustomLayouts.Count; $i++) {
    $cl = $m.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $NEW_DATE
        }
    }
}

# --- 2) & 3) Per-slide background + accent rectangle recolor ---
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    $s.Background.Fill.ForeColor.RGB = $ACCENT_BGR

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.Name -eq "Rectangle 3" -or $sh.Name -eq "Rectangle 5") {
            $sh.Fill.ForeColor.RGB = $GOLD_BGR
        }
    }
}

# --- 4) Slide 6 permissions table recolor + tiny relayout height fix ---
$s6 = $p.Slides.Item(6)
for ($j = 1; $j -le $s6.Shapes.Count; $j++) {
    $sh = $s6.Shapes.Item($j)
    if ($sh.HasTable) {
        $tbl = $sh.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cell = $tbl.Cell($r, $c)
                $cell.Shape.Fill.ForeColor.RGB = $ACCENT_BGR
            }
        }
        $sh.Height = 401.8722834645669
    }
}
